# Updated cryptos list on Wed Jun 12 09:58:27 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.653.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.91%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.531.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.48%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.07%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.529.66"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.47%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.04%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.35%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.50%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +3.59%  "

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.427"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.36%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -0.08%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.42%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.120.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "

# Row 16 - WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.525.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.21%  "

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.499.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.62%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -0.01%  "

# Row 19 - Polkadot
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.56%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -0.28%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "447.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.94%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.81%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  -2.69%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.66%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  +8.09%  "

# Row 26 - WrappedeETH
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.669.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.23%  "

# Row 27 / Row 28 - coins swapped order: InternetComputer(DFINITY) now ranks above Dai
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.95%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.18%  "

# Row 29 - RenderToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.17%  "

# Row 31 - Fetch.AI
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.16%  "

# Row 32 - Binance-PegBSC-USD
$ws.Range("E32").Value = "  -0.06%  "

# Row 33 - Kaspa
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.166"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.68%  "

# Row 34 - EthereumClassic
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.30%  "

# Row 35 - NEARProtocol
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.13%  "

# Row 36 - RenzoRestakedETH
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.518.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.23%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -2.47%  "

# Row 38 - Aptos
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.71%  "

# Row 39 - USDe
$ws.Range("E39").Value = "  +0.03%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  -0.10%  "

# Row 41 - Monero
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "177.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.48%  "

# Row 42 - Hedera
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0889"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.27%  "

# Row 43 - Stacks
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.79%  "

# Row 44 - Filecoin
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.73%  "

# Row 45 - Mantle
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.885"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.71%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.59%  "

# Row 47 - OKB
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.86%  "

# Row 48 - dogwifhat
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.12%  "

# Row 49 - ONDO
$ws.Range("E49").Value = "  +3.65%  "

# Row 50 - Cosmos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.42%  "

# Row 51 - SuiNetwork
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.73%  "

